$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("grouped matches")

$updates = @(
    @{Row=2; B='P8211'; C='{''eft:vidyakaraprabha''}'},
    @{Row=3; B='P4CZ15137'; C='{''eft:kumarakalasa''}'},
    @{Row=4; B='P4263'; C='{''eft:dge-ba-dpal''}'},
    @{Row=5; B='P4258'; C='{''eft:dpal-byor''}'},
    @{Row=6; B='P0TMP098'; C='{''eft:jinavara''}'},
    @{Row=7; B='https://lod.dila.edu.tw/resource.php?id=A000089'; C='{''eft:siladharma''}'},
    @{Row=8; B='P8245'; C='{''eft:buddhakaravarma''}'},
    @{Row=9; B='P8222'; C='{''eft:jnanasiddhi'', ''eft:jnanasidhi''}'},
    @{Row=10; B='P4CZ16780'; C='{''eft:manjusrigarbha''}'},
    @{Row=11; B='P8267'; C='{''eft:vijayasila''}'},
    @{Row=12; B='P00KG07267'; C='{''eft:sarvanyadeva'', ''eft:sarvajnadeva''}'},
    @{Row=13; B='P8220'; C='{''eft:devacandra''}'},
    @{Row=14; B='P0TMP080'; C='{''eft:hwa-shang-zab-mo''}'},
    @{Row=15; B='P8261'; C='{''eft:munivarman'', ''eft:munivarma''}'},
    @{Row=16; B='P4255'; C='{''eft:t-jnanagarbha'', ''eft:yesh-nyingpo'', ''eft:ye-shes-snying-po''}'},
    @{Row=17; B='P8268'; C='{''eft:buddhaprabha''}'},
    @{Row=18; B='P3456'; C='{''eft:tshul-khrims-rgyal-ba''}'},
    @{Row=19; B='P8249'; C='{''eft:dharmakara''}'},
    @{Row=20; B='P0TMPT007'; C='{''eft:rnam-par-mi-rtog-pa''}'},
    @{Row=21; B='P5651'; C='{''eft:pa-tshab-nyi-ma-grags''}'},
    @{Row=22; B='P4CZ16819'; C='{''eft:sakyaprabha''}'},
    @{Row=23; B='P753'; C='{''eft:rin-chen-bzang-po''}'},
    @{Row=24; B='P8093'; C='{''eft:kamalagupta''}'},
    @{Row=25; B='P8213'; C='{''eft:t-vidyakarasimha'', ''eft:vidyakarasimha''}'},
    @{Row=26; B='P8171'; C='{''eft:dharmasribhadra''}'},
    @{Row=27; B='P3379'; C='{''eft:dipamkara-srijnana'', ''eft:dipamkarasrijnana''}'},
    @{Row=28; B='P8205'; C='{''eft:yesh-d-ye-shes-sde-'', ''eft:zhang-yesh-d-'', ''eft:ye-shes-sde'', ''eft:yesh-d-'', ''eft:band-yesh-de'', ''eft:band-yesh-d-''}'},
    @{Row=29; B='P8182'; C='{''eft:ska-ba-dpal-brtsegs'', ''eft:dpal-brtsegs'', ''eft:ban-de-dpal-brtsegs'', ''eft:kawa-paltsek-under-the-name-paltsek-raksita-'', ''eft:paltsek''}'},
    @{Row=30; B='P8206'; C='{''eft:celu''}'},
    @{Row=31; B='P0TMP092'; C='{''eft:anandasri-s-''}'},
    @{Row=32; B='P8269'; C='{''eft:dgon-gling-rma''}'},
    @{Row=33; B='P2548'; C='{''eft:prajnavarma'', ''eft:prajnavarman''}'},
    @{Row=34; B='P3214'; C='{''eft:danasila''}'},
    @{Row=35; B='P8219'; C='{''eft:visuddhasimha''}'},
    @{Row=36; B='P8183'; C='{''eft:klu-i-rgyal-mtshan'', ''eft:cog-ro-klu-i-rgyal-mtshan''}'},
    @{Row=37; B='P8217'; C='{''eft:t-jnanagarbha'', ''eft:jnanagarbha''}'},
    @{Row=38; B='P3709'; C='{''eft:phakpa-sherab''}'},
    @{Row=39; B='P3285'; C='{''eft:sakya-yesh-''}'},
    @{Row=40; B='P1KG8854'; C='{''eft:surendrabodhi'', ''eft:silendrabodhi'', ''eft:srilendrabodhi''}'},
    @{Row=41; B='P8263'; C='{''eft:leki-d-''}'},
    @{Row=42; B='P4242'; C='{''eft:sherab-lekpa''}'},
    @{Row=43; B='P8265'; C='{''eft:ratnaraksita''}'},
    @{Row=44; B='P0RK8'; C='{''eft:dharmapala''}'},
    @{Row=45; B='P2637'; C='{''eft:trakpa-gyaltsen''}'},
    @{Row=46; B='P8273'; C='{''eft:rinchen-tso'', ''eft:rin-chen-tsho''}'},
    @{Row=47; B='P4259'; C='{''eft:ban-de-dpal-gyi-lhun-po'', ''eft:dpal-gyi-lhun-po'', ''eft:palgyi-lh-npo''}'},
    @{Row=48; B='P8209'; C='{''eft:jinamitra'', ''eft:jinamitra-k-'', ''eft:dzi-na-mi-tra-k-''}'},
    @{Row=49; B='P2956'; C='{''eft:krsnapandita''}'},
    @{Row=50; B='P8151'; C='{''eft:gayadhara''}'},
    @{Row=51; B='P0TMP104'; C='{''eft:punyasambhava''}'},
    @{Row=52; B='P8228'; C='{''eft:surendrabodhi''}'},
    @{Row=53; B='P8266'; C='{''eft:ch-nyi-tsultrim'', ''eft:dharmatasila''}'},
    @{Row=54; B='?'; C='{''eft:sakyasena''}'},
    @{Row=55; B='P8260'; C='{''eft:dpal-dbyangs''}'}
)

foreach ($u in $updates) {
    $ws.Cells.Item($u.Row, 2).Value = $u.B
    $ws.Cells.Item($u.Row, 3).Value = $u.C
}
